$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Property")

# Fill in server data row (row 2)
$ws.Range("A2").Value = "SqlServer_1"
$ws.Range("C2").Value = "127.0.0.1"
$ws.Range("B2").Value = "000107001"
$ws.Range("D2").Value = 7001
$ws.Range("E2").Value = 123456

# Update the active selection on the sheet
$ws.Range("E4").Select()
